$d = $word.ActiveDocument

# --- 1. Heading text: "Effect" -> "the Effects" ---
$d.Content.Find.Execute(
    "Visualizing Effect of Data Transformations on Errors",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Visualizing the Effects of Data Transformations on Errors", 2)

# --- 2. Rename the heading bookmark to match the new title/slug ---
$oldBookmarkName = "visualizing-effect-of-data-transformations-on-errors"
$newBookmarkName = "visualizing-the-effects-of-data-transformations-on-errors"
$bm = $d.Bookmarks($oldBookmarkName)
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add($newBookmarkName, $bmRange)

# --- 3. Rewrite the abstract body paragraph ---
$oldBody = "In many omics analyses, a primary step is data transformation. A transformation is generally employed to convert proportional error (variance) to additive error, which most statistical methods appropriately handle. However, omics data frequently contain error sources that result in both additive and proportional errors. To our knowledge, there has not been a systematic study on detecting proportional error in omics data, or the effect of transformations on the error structure. In this work we demonstrate that a set of three simple graphs can facilitate the detection of proportional and mixed error in omics data. The graphs plot the 1) absolute range, 2) standard deviation and 3) relative standard deviation across replicates against the mean signal across replicates. In addition to showing the presence of different types of error, they also readily demonstrate the effect of various transformations on the error structure as well. Using these graphical summaries we find that the log-transform is the most effective method of the common methods employed for removing proportional error."
$newBody = "In many omics analyses, a primary step in the data analysis is the application of a transformation to the data. Transformations are generally employed to convert proportional error (variance) to additive error, which most statistical methods appropriately handle. However, omics data frequently contain error sources that result in both additive and proportional errors. To our knowledge, there has not been a systematic study on detecting the presence of proportional error in omics data, or the effect of transformations on the error structure. In this work, we demonstrate a set of three simple graphs which facilitate the detection of proportional and mixed error in omics data when multiple replicates are available. The three graphs illustrate proportional and mixed error in a visually compelling manner that is both straight-forward to recognize and to communicate. The graphs plot the 1) absolute range, 2) standard deviation and 3) relative standard deviation against the mean signal across replicates. In addition to showing the presence of different types of error, these graphs readily demonstrate the effect of various transformations on the error structure as well. Using these graphical summaries we find that the log-transform is the most effective method of the common methods employed for removing proportional error."

$d.Content.Find.Execute($oldBody, $true, $false, $false, $false, $false, $true, 1, $false, $newBody, 2)
